$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6, columns A-J
$data = @(
    @(1, 9, 2, 6, 5, -3, 3, 34, 5, "train_dim2_1"),
    @(2, 5, 0, 4, 5, -1, 5, 56, 5, "train_dim2_1"),
    @(3, 8, 1, 3, 2, -5, 1, 12, 5, "train_dim2_1"),
    @(4, 5, 2, 3, 6, -2, 4, 45, 5, "train_dim2_1"),
    @(5, 9, 1, 5, 3, -4, 2, 23, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

$ws.Range("I1").Select()
